$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand new "2022-Q3" sheet right after "总计", and populate it
#    with the fund-holding detail rows (same layout as the other quarterly
#    detail sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名)
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Columns B (fund code, has leading zeros) and D,E,F,G (stored as plain text,
# e.g. "13.72") are text in this workbook, not numbers -- force text format
# before assigning so Excel doesn't coerce them into numeric cells.
$q3.Range("B2:B14").NumberFormat = "@"
$q3.Range("D2:G14").NumberFormat = "@"

$q3Rows = @(
    @(0,  "515450", "南方标普中国A股大盘红利低波50ETF",          "2.17", "99.66", "4.04", "0.0877", 4),
    @(1,  "012708", "东方红中证东方红红利低波动指数A",            "3.27", "93.80", "1.66", "0.0543", 4),
    @(2,  "008115", "天弘中证红利低波动100指数C",                 "2.44", "94.56", "1.97", "0.0481", 5),
    @(3,  "080005", "长盛量化红利混合",                           "1.89", "61.68", "2.33", "0.0440", 8),
    @(4,  "008114", "天弘中证红利低波动100指数A",                 "1.89", "94.56", "1.97", "0.0372", 5),
    @(5,  "515100", "景顺长城中证红利低波动100ETF",               "1.62", "98.63", "2.06", "0.0334", 5),
    @(6,  "512190", "浙商汇金中证浙江凤凰行动50ETF",              "0.48", "98.92", "4.66", "0.0224", 5),
    @(7,  "012709", "东方红中证东方红红利低波动指数C",            "0.67", "93.80", "1.66", "0.0111", 4),
    @(8,  "007751", "景顺长城中证沪港深红利成长低波动指数A",      "0.67", "90.27", "1.51", "0.0101", 8),
    @(9,  "009384", "摩根士丹利华鑫MSCI中国A股指数增强A",        "0.39", "90.98", "1.07", "0.0042", 7),
    @(10, "005126", "银河量化稳进混合",                           "0.13", "55.69", "1.09", "0.0014", 8),
    @(11, "007760", "景顺长城中证沪港深红利成长低波动指数C",      "0.06", "90.27", "1.51", "0.0009", 8),
    @(12, "014866", "摩根士丹利华鑫MSCI中国A股指数增强C",        "0.00", "90.98", "1.07", $null,   7)
)

for ($i = 0; $i -lt $q3Rows.Count; $i++) {
    $r = $i + 2
    $row = $q3Rows[$i]
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").Value = $row[5]
    if ($null -ne $row[6]) {
        $q3.Range("G$r").Value = $row[6]
    } else {
        # last fund's "持有市值" is a genuine number 0, not text, in this sheet
        $q3.Range("G$r").NumberFormat = "General"
        $q3.Range("G$r").Value = 0
    }
    $q3.Range("H$r").Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. Insert a new row into "总计" summarising 2022-Q3, pushing the existing
#    quarters down by one row. The A column is a plain 0-based row index, so
#    every pre-existing row's index needs to be bumped by one too.
# ---------------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

for ($r = 8; $r -ge 3; $r--) {
    $zongji.Range("A$r").Value = $r - 2
}

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 13
$zongji.Range("D2").Value = 0.35
